# Swap the order of "Recorded By" entries in column G so that "System"
# is listed first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# Only applies to the two specific author emails that were re-ordered upstream
# (dnasr281@gmail.com and admin@admin.com); "backup@backdoor.com, System" is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq "dnasr281@gmail.com, System") {
        $cell.Value2 = "System, dnasr281@gmail.com"
    }
    elseif ($val -eq "admin@admin.com, System") {
        $cell.Value2 = "System, admin@admin.com"
    }
}
